# Apply updated coin price / label values to the "cryptos" worksheet.
# All target cells are plain text cells (e.g. "244.20"), so we force the
# NumberFormat to Text ("@") before assigning the new value. This prevents
# Excel from auto-converting numeric-looking strings into numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        [string]$Address,
        [string]$Value
    )
    $cell = $ws.Range($Address)
    $cell.NumberFormat = "@"
    $cell.Value = $Value
}

# Column D (Price) updates
Set-TextValue "D2"  "244.30"
Set-TextValue "D3"  "23.04"
Set-TextValue "D4"  "5.415"
Set-TextValue "D5"  "0.05973"
Set-TextValue "D6"  "3.464"
Set-TextValue "D7"  "6.514"
Set-TextValue "D8"  "0.8133"
Set-TextValue "D9"  "0.9133"
Set-TextValue "D10" "0.1410"
Set-TextValue "D11" "0.07415"
Set-TextValue "D12" "0.03236"
Set-TextValue "D13" "0.03091"
Set-TextValue "D14" "0.09363"
Set-TextValue "D15" "3.848"
Set-TextValue "D16" "0.001558"
Set-TextValue "D17" "0.04676"
Set-TextValue "D18" "0.0005941"
Set-TextValue "D19" "0.006093"

# Row 20: price + label change
Set-TextValue "D20" "0.005000"
Set-TextValue "E20" "19HotbitTokenHTB"

Set-TextValue "D21" "0.0009830"
Set-TextValue "D22" "0.00007800"

Set-TextValue "D25" "0.3206"
Set-TextValue "D26" "0.1303"
Set-TextValue "D27" "0.0002900"

Set-TextValue "D40" "0.03931"

# Row 41: price + label change
Set-TextValue "D41" "0.006235"
Set-TextValue "E41" "40KickTokenKICKBestin24h"

Set-TextValue "D42" "0.1075"
Set-TextValue "D43" "0.003000"
Set-TextValue "D44" "0.006552"
Set-TextValue "D45" "0.00005243"

Set-TextValue "D48" "0.8701"
Set-TextValue "D49" "0.002284"
Set-TextValue "D50" "0.00002100"
Set-TextValue "D51" "0.0002000"
